$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header figures: updated totals for this statement of account ---
$ws.Range("E11").Value = 690080
$ws.Range("F13").Value = 17

# --- Trim the periods-in-arrears table (remove settled periods) ---
# Delete from the bottom up so row numbers of not-yet-deleted rows stay valid.
# Removes MARIA PAULA's 2503/2502 rows (31:32) and CAMILA's 2505/2504 rows (34:35),
# leaving one row per remaining open period and shifting STEFFANY's closing row up.
$ws.Rows("34:35").Delete()
$ws.Rows("31:32").Delete()

# --- Update the values on the rows that remain ---
# MARIA PAULA ARENAS JIMENEZ - 2504 : valor mora updated
$ws.Range("G30").Value = 1423500

# CAMILA FUENTES RODRIGUEZ - 2506 : valor mora updated
$ws.Range("G31").Value = 1423500

# STEFFANY ANDREA MEZA RENGIFO - new period 2508 added as the latest period in arrears
$ws.Range("E32").Value = "2508"
$ws.Range("F32").Value = 56000
